# structure revamp for future tools and internationalization
#
# Two blog rows ("What is an Aha Moment" and "13 cognitive biases to avoid
# in daily life") are promoted to the top of the table (rows 2 and 3). All
# the rows that used to sit between the old top and the old "Aha Moment"
# row cascade down by one. The "updated" timestamps (column D) for the two
# promoted rows are refreshed; their "published" timestamps (column G) -
# and every other value for the cascaded rows - travel with the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for rows 2..18 (row 1 is the header, rows 19..21 are
# untouched) expressed as A/B/C/D/F/G tuples.
$rows = @(
    @{ A = "What is an Aha Moment";                          B = "what-is-an-aha-moment";                                      C = "What Is An Aha Moment";                          D = 44983.66443622751; G = 44968.82335802943 },
    @{ A = "13 cognitive biases to avoid in daily life";      B = "cognitive-biases-to-avoid-in-daily-life";                    C = "13 Cognitive Biases To Avoid In Daily Life";      D = 44983.66333873353; G = 44969.74406943242 },
    @{ A = "Key setup for problem solving";                   B = "key-setup-for-problem-solving";                              C = "Key Setup For Problem Solving";                   D = 44970.84635276953; G = 44969.9128825571  },
    @{ A = "What is the 80 percent rule";                     B = "what-is-the-80-percent-rule";                                C = "What Is The 80 Percent Rule";                     D = 44969.53807938268; G = 44968.82652412132 },
    @{ A = "What is K Level Thinking";                        B = "what-is-k-level-thinking";                                   C = "What Is K Level Thinking";                        D = 44969.5374997161;  G = 44968.81650771474 },
    @{ A = "What is the ambiguity effect";                    B = "what-is-the-ambiguity-effect";                               C = "What Is The Ambiguity Effect";                    D = 44969.53642034208; G = 44968.83411760339 },
    @{ A = "What is the Skyscraper Method";                   B = "what-is-the-skyscraper-method";                              C = "What Is The Skyscraper Method";                   D = 44969.53607302732; G = 44968.82335802943 },
    @{ A = "What is the Technology Window";                   B = "what-is-the-technology-window";                              C = "What Is The Technology Window";                   D = 44969.53543811679; G = 44968.82226430001 },
    @{ A = "What is the single source of truth principle";    B = "what-is-the-single-source-of-truth-principle";               C = "What Is The Single Source Of Truth Principle";    D = 44969.53487175713; G = 44968.8283912778  },
    @{ A = "What is the self-fulfilling prophecy";            B = "what-is-the-self-fulfilling-prophecy";                       C = "What Is The Self-Fulfilling Prophecy";            D = 44969.53473463681; G = 44968.83482908575 },
    @{ A = "What is the quantitative accumulation effect";    B = "what-is-the-quantitative-accumulation-effect";               C = "What Is The Quantitative Accumulation Effect";    D = 44969.53425319671; G = 44968.82762568397 },
    @{ A = "What is the Narrow Path Principle";                B = "what-is-the-narrow-path-principle";                          C = "What Is The Narrow Path Principle";               D = 44969.5339683831;  G = 44968.836970087   },
    @{ A = "What is the Dunning-Krueger-Effect";               B = "what-is-the-dunning-krueger-effect";                         C = "What Is The Dunning-Krueger-Effect";              D = 44969.53324986236; G = 44968.83607997302 },
    @{ A = "What is the believe bias";                         B = "what-is-the-believe-bias";                                   C = "What Is The Believe Bias";                        D = 44969.53235714744; G = 44968.83201703853 },
    @{ A = "What is the Cheerleader effect";                   B = "what-is-the-cheerleader-effect";                             C = "What Is The Cheerleader Effect";                  D = 44969.53145739079; G = 44968.83272623207 },
    @{ A = "What is the Bystander Effect";                     B = "what-is-the-bystander-effect";                               C = "What Is The Bystander Effect";                    D = 44969.5313018782;  G = 44968.83799102488 },
    @{ A = "What is the barnum effect";                        B = "what-is-the-barnum-effect";                                  C = "What Is The Barnum Effect";                       D = 44969.53007293754; G = 44968.83331997384 }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A   # A: document/h1-source title
    $ws.Cells.Item($r, 2).Value = $row.B   # B: url slug
    $ws.Cells.Item($r, 3).Value = $row.C   # C: h1 (title case)
    $ws.Cells.Item($r, 4).Value = $row.D   # D: updated
    $ws.Cells.Item($r, 6).Value = $row.A   # F: document (mirrors A)
    $ws.Cells.Item($r, 7).Value = $row.G   # G: published
    $r++
}
